$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New presence date column for 24.02, mirroring the existing 23.02 (D) column
# Force E3 to text so the numeric-looking "24.02" is stored as a string, not a date/number
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "24.02"
$ws.Range("E4").Value = "x"
$ws.Range("E5").Value = "x"
$ws.Range("E6").Value = "x"
$ws.Range("E7").Value = "x"
$ws.Range("E8").Value = "x"

# Match the style used by the existing D column (plain/default style)
$ws.Range("D3:D8").Copy()
$ws.Range("E3:E8").PasteSpecial(-4122)

# Final cursor position left by the author
$ws.Range("G7").Select()
